$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure edited numeric-looking strings remain text (matching original inline string cells)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "305.71"
$ws.Range("E2").Value = "-0.77%"
$ws.Range("D3").Value = "36.39"
$ws.Range("E3").Value = "-0.75%"
$ws.Range("D4").Value = "5.037"
$ws.Range("E4").Value = "-0.10%"
$ws.Range("D5").Value = "0.07935"
$ws.Range("E5").Value = "0.79%"
$ws.Range("D6").Value = "2.125"
$ws.Range("E6").Value = "-2.75%"
$ws.Range("D7").Value = "7.968"
$ws.Range("E7").Value = "-1.04%"
$ws.Range("D8").Value = "0.9255"
$ws.Range("E8").Value = "-0.13%"
$ws.Range("D9").Value = "0.09806"
$ws.Range("E9").Value = "-1.37%"
$ws.Range("D10").Value = "0.1864"
$ws.Range("E10").Value = "-0.92%"
$ws.Range("D11").Value = "0.09001"
$ws.Range("E11").Value = "3.71%"
$ws.Range("D12").Value = "0.03602"
$ws.Range("E12").Value = "0.14%"
$ws.Range("D13").Value = "0.09932"
$ws.Range("E13").Value = "-0.10%"
$ws.Range("D14").Value = "0.001440"
$ws.Range("E14").Value = "-2.32%"
$ws.Range("D15").Value = "0.005636"
$ws.Range("E15").Value = "-0.50%"
$ws.Range("D16").Value = "3.481"
$ws.Range("E16").Value = "0.58%"
$ws.Range("D17").Value = "4.140"
$ws.Range("E17").Value = "1.94%"
$ws.Range("E18").Value = "7.66%"
$ws.Range("D19").Value = "0.3423"
$ws.Range("E19").Value = "-0.34%"
$ws.Range("D20").Value = "0.1326"
$ws.Range("E20").Value = "-0.19%"
$ws.Range("D21").Value = "5.172"
$ws.Range("E21").Value = "4.97%"
$ws.Range("D22").Value = "0.2246"
$ws.Range("E22").Value = "2.05%"
$ws.Range("D23").Value = "0.04580"
$ws.Range("E23").Value = "-0.85%"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").Value = "-0.74%"
$ws.Range("D25").Value = "0.004824"
$ws.Range("E25").Value = "-7.40%"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").Value = "-7.24%"
$ws.Range("D27").Value = "0.0004738"
$ws.Range("E27").Value = "74.29%"
$ws.Range("D39").Value = "0.01878"
$ws.Range("E39").Value = "2.48%"
$ws.Range("D40").Value = "0.04927"
$ws.Range("E40").Value = "3.45%"
$ws.Range("D41").Value = "0.007755"
$ws.Range("E41").Value = "-2.19%"
$ws.Range("D42").Value = "0.1399"
$ws.Range("E42").Value = "-0.90%"
$ws.Range("D43").Value = "0.007750"
$ws.Range("E43").Value = "1.99%"
$ws.Range("D44").Value = "0.002109"
$ws.Range("E44").Value = "-1.09%"
$ws.Range("D45").Value = "0.01124"
$ws.Range("E45").Value = "11.12%"
$ws.Range("D46").Value = "0.00006427"
$ws.Range("E46").Value = "1.66%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.12%"
$ws.Range("D49").Value = "51.72"
$ws.Range("E49").Value = "42.77%"
$ws.Range("D50").Value = "0.001899"
$ws.Range("E50").Value = "-29.43%"
$ws.Range("D51").Value = "0.00002099"
$ws.Range("E51").Value = "-0.12%"
